# Apply the freeCrm keyword-driven test framework update.
#
# Summary of changes:
#  - "signup" sheet becomes the active/selected sheet (was "login").
#  - "login" sheet keeps selection on E8 (was C8) and is no longer tabSelected.
#  - "login" sheet E2 browser stays "chrome" (unchanged), but the "get/verify home
#    page header" rows are updated to target freeCRM's header markup instead of
#    HubSpot's Sales Dashboard.
#  - "signup" sheet E2 browser changes from "firefox" to "chrome", and the login
#    URL changes from HubSpot to the freeCRM classic login/launch URL, including
#    the underlying hyperlink target.

$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("login")
$wsSignup = $wb.Worksheets.Item("signup")

# ---- login sheet updates ----

# Row 7: "verify home page header"
$wsLogin.Range("C7").Value = "//td[@class='headertext']"
$wsLogin.Range("D7").Value = "getText"
$wsLogin.Range("E7").Value = "User: Mehraj Ismayilov"

# Row 8: "get home page header text"
$wsLogin.Range("C8").Value = "headertext"
$wsLogin.Range("E8").Value = "User: Mehraj Ismayilov"

# ---- signup sheet updates ----

# Row 2: browser changes from firefox to chrome
$wsSignup.Range("E2").Value = "chrome"

# Row 3: launch url changes from HubSpot login to freeCRM classic url,
# including the hyperlink target behind the cell.
$wsSignup.Range("E3").Value = "https://classic.freecrm.com"
foreach ($hl in $wsSignup.Hyperlinks) {
    $hl.Address = "https://classic.freecrm.com"
}

# ---- selection / active sheet updates ----

# "login" selection moves to E8, but it is no longer the tab-selected sheet.
$wsLogin.Range("E8").Select()

# "signup" becomes the active / tab-selected sheet.
$wsSignup.Activate()
